# "changed xpath for testcases"
#
# Summary of edits applied to Webdata_TestData.xlsx:
#  1. AddCustomerChild (sheet19): append two new columns (Q1, R1) of test data,
#     and the active/selected cell moves off this sheet.
#  2. GenerateInvoice (sheet21): the row is re-shaped - several leading
#     columns collapse into a single "Ashish" value, and the "abc" test
#     placeholder is replaced by a real description, shifting everything
#     after it left by a few columns.
#  3. Reports (sheet27): becomes the newly active sheet and gains a new
#     trailing column (Q1) with a Commission Product reference.
#  4. The now-unused shared string "abc" disappears automatically once
#     nothing references it any more.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. AddCustomerChild -- add two new trailing cells
# ---------------------------------------------------------------------
$wsChild = $wb.Worksheets.Item("AddCustomerChild")
$wsChild.Range("Q1").Value = "22 - Brian Smith"
$wsChild.Range("R1").Value = "Customer 22 receives invoices. "

# ---------------------------------------------------------------------
# 2. GenerateInvoice -- reshape row 1
#    old: A..S = admin, Webdata@123, Web Data US, Web Data US,
#                Direct Customer, Customer A, Monthly, 1, Monthly,
#                pre paid, 10/17/2016, abc, 10/17/2016, 4, Test Code,
#                Test Code1, TC-01, 5, 10
#    new: A..O = admin, Webdata@123, Web Data US, Ashish, Monthly,
#                pre paid, 10/17/2016, Product Code1 Description,
#                10/17/2016, 4, Test Code, Test Code1, TC-01, 5, 10
# ---------------------------------------------------------------------
$wsInvoice = $wb.Worksheets.Item("GenerateInvoice")

# Stage the values that need to shift left, in far-away columns first so
# that copying back does not clobber a value before it has been staged.
$wsInvoice.Range("G1").Copy($wsInvoice.Range("AA1"))
$wsInvoice.Range("J1").Copy($wsInvoice.Range("AB1"))
$wsInvoice.Range("K1").Copy($wsInvoice.Range("AC1"))
$wsInvoice.Range("M1").Copy($wsInvoice.Range("AD1"))
$wsInvoice.Range("N1").Copy($wsInvoice.Range("AE1"))
$wsInvoice.Range("O1").Copy($wsInvoice.Range("AF1"))
$wsInvoice.Range("P1").Copy($wsInvoice.Range("AG1"))
$wsInvoice.Range("Q1").Copy($wsInvoice.Range("AH1"))
$wsInvoice.Range("R1").Copy($wsInvoice.Range("AI1"))
$wsInvoice.Range("S1").Copy($wsInvoice.Range("AJ1"))

# Copy the staged values into their final, shifted positions.
$wsInvoice.Range("AA1").Copy($wsInvoice.Range("E1"))
$wsInvoice.Range("AB1").Copy($wsInvoice.Range("F1"))
$wsInvoice.Range("AC1").Copy($wsInvoice.Range("G1"))
$wsInvoice.Range("AD1").Copy($wsInvoice.Range("I1"))
$wsInvoice.Range("AE1").Copy($wsInvoice.Range("J1"))
$wsInvoice.Range("AF1").Copy($wsInvoice.Range("K1"))
$wsInvoice.Range("AG1").Copy($wsInvoice.Range("L1"))
$wsInvoice.Range("AH1").Copy($wsInvoice.Range("M1"))
$wsInvoice.Range("AI1").Copy($wsInvoice.Range("N1"))
$wsInvoice.Range("AJ1").Copy($wsInvoice.Range("O1"))

# New literal text values.
$wsInvoice.Range("D1").Value = "Ashish"
$wsInvoice.Range("H1").Value = "Product Code1 Description"

# Remove the old trailing cells (including the staging area).
$wsInvoice.Range("P1:AJ1").ClearContents()

# ---------------------------------------------------------------------
# 3. Reports -- add a new trailing cell
# ---------------------------------------------------------------------
$wsReports = $wb.Worksheets.Item("Reports")
$wsReports.Range("Q1").Value = "Commission Product"

# ---------------------------------------------------------------------
# 4. Selection / active sheet bookkeeping
# ---------------------------------------------------------------------
$null = $wsChild.Activate()
$null = $wsChild.Range("R1").Select()

$null = $wsInvoice.Activate()
$null = $wsInvoice.Range("L1").Select()

$null = $wsReports.Activate()
$null = $wsReports.Range("Q1").Select()
